$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36/37 swap: Kaspa <-> Monero (Monero now ranked above Kaspa) ---
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'157.57"
$ws.Range("E36").Value = "  +2.99%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.115"
$ws.Range("E37").Value = "  -3.98%  "

# --- Price (D) / Volume 1h (E) refresh for all other rows ---
$ws.Range("D2").Value = "69.602.05"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.501.47"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'574.38"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "'166.50"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").Value = "2.500.03"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "2.958.34"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "69.562.71"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'24.67"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "2.508.32"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").Value = "'7.43"
$ws.Range("E20").Value = "  -4.72%  "
$ws.Range("D21").Value = "'348.10"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "'1.94"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'70.67"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("D26").Value = "'3.92"
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("D28").Value = "2.630.05"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "'457.57"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("E33").Value = "  -5.26%  "
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D38").Value = "'19.05"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "'38.02"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("E46").Value = "  -7.18%  "
$ws.Range("D47").Value = "'141.11"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("E49").Value = "  -2.83%  "
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "'0.577"
$ws.Range("E51").Value = "  -0.65%  "
